$d = $word.ActiveDocument

# Note: this runtime's Find.Execute replaces against the whole document
# regardless of which Range/Find object it is invoked on, so we use
# wdReplaceOne (1) and apply the replacements strictly in document order.
# That way each call only touches the next remaining (first) occurrence
# of the search text, correctly handling duplicate strings such as the
# two "63÷8=" cells.

$replacements = @(
    @{ Old = "2024-11-21 Thursday"; New = "2024-11-22 Friday" },

    @{ Old = "63÷8="; New = "98÷4=" },
    @{ Old = "90÷9="; New = "33÷3=" },
    @{ Old = "97÷2="; New = "15÷9=" },
    @{ Old = "64÷3="; New = "16÷2=" },
    @{ Old = "23÷3="; New = "80÷2=" },

    @{ Old = "45÷5="; New = "80÷3=" },
    @{ Old = "21÷9="; New = "82÷9=" },
    @{ Old = "41÷2="; New = "93÷5=" },
    @{ Old = "32÷5="; New = "77÷5=" },
    @{ Old = "49÷9="; New = "67÷8=" },

    @{ Old = "83÷7="; New = "27÷3=" },
    @{ Old = "15÷2="; New = "78÷7=" },
    @{ Old = "63÷8="; New = "29÷3=" },
    @{ Old = "14÷2="; New = "67÷8=" },
    @{ Old = "66÷6="; New = "15÷8=" },

    @{ Old = "39÷2="; New = "96÷5=" },
    @{ Old = "47÷8="; New = "35÷4=" },
    @{ Old = "69÷7="; New = "19÷7=" },
    @{ Old = "48÷5="; New = "45÷9=" },
    @{ Old = "66÷3="; New = "77÷3=" },

    @{ Old = "12÷2="; New = "42÷2=" },
    @{ Old = "69÷4="; New = "37÷2=" },
    @{ Old = "63÷7="; New = "28÷3=" },
    @{ Old = "53÷2="; New = "90÷2=" },
    @{ Old = "43÷5="; New = "60÷9=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 1) | Out-Null
}
